# Commit com compressão de video
# Mark all remaining open backlog items ("Realizado") as done (TRUE),
# which ripples into the burndown totals/averages and the chart cache.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("I12:I25").Value = $true
$ws.Range("I28:I37").Value = $true

$ws.Range("I18").Select()
